$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# NOTE: the COM ColumnWidth setter here stores width + 5/6 character units
# (mirrors Excel's default-font pixel rounding), so subtract that offset to
# land on the exact target "characters" width used by the target file.
$pad = 0.83333333333333

$ws.Columns.Item(2).ColumnWidth = 30 - $pad
$ws.Columns.Item(5).ColumnWidth = 15 - $pad
$ws.Columns.Item(11).ColumnWidth = 40 - $pad

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 25
$ws.Rows.Item(2).RowHeight = 30

# --- Fix vendor name text (drop comma after WESTSIDE) ---
$ws.Range("B2").Value = "WESTSIDE Sjr Zion, Survey"

# --- Alignment / formatting for row 2 ---
# A2: S.No. -> centered horizontally & vertically
$ws.Range("A2").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A2").VerticalAlignment = -4108    # xlCenter

# B2 (Vendor name) and K2 (HSN Codes): vertical center + wrap text
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").WrapText = $true
$ws.Range("K2").VerticalAlignment = -4108
$ws.Range("K2").WrapText = $true

# C2:J2 : vertical center only
$ws.Range("C2:J2").VerticalAlignment = -4108
